# Remove six account rows from the "Export" sheet's data table.
# Accounts removed (by row, in the original layout):
#   row 4  -> 004332747 / LOHRAN   / 62551.62
#   row 6  -> 004468717 / HELOISA  / 16460.95
#   row 7  -> 004363260 / LARISSA  / 14486.98
#   row 8  -> 003553997 / MIRELLA  / 10448.06
#   row 9  -> 004565108 / GUSTAVO  / 5000
#   row 11 -> 004228090 / GUSTAVO  / 2550.57
#
# Deleting whole rows (EntireRow / Rows(n).Delete) shifts everything below
# up, matching the diff where the remaining rows simply close the gaps.
# Rows are deleted from the bottom up so earlier row numbers stay valid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(11, 9, 8, 7, 6, 4)

foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}
